# Update the cryptos worksheet with refreshed Price (D) and Volume(1h) (E)
# values, as produced by the scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  D = '79.192.96';   E = '  +3.83%  ' }
    @{ Row = 3;  D = '3.143.67';    E = '  +2.00%  ' }
    @{ Row = 4;  D = '1.00';        E = '  +0.01%  ' }
    @{ Row = 5;  D = '203.99';      E = '  +2.87%  ' }
    @{ Row = 6;  D = '621.66';      E = '  +0.49%  ' }
    @{ Row = 7;  D = '0.262';       E = '  +25.21%  ' }
    @{ Row = 8;  D = '0.999';       E = '  -0.03%  ' }
    @{ Row = 9;  D = '0.585';       E = '  +6.09%  ' }
    @{ Row = 10; D = '3.145.39';    E = '  +2.11%  ' }
    @{ Row = 11; D = '0.579';       E = '  +31.04%  ' }
    @{ Row = 12; D = '0.0000249';   E = '  +27.96%  ' }
    @{ Row = 13; D = '0.164';       E = '  +1.82%  ' }
    @{ Row = 14; D = '3.722.30';    E = '  +1.99%  ' }
    @{ Row = 15; D = '5.23';        E = '  -0.33%  ' }
    @{ Row = 16; D = '31.00';       E = '  +5.71%  ' }
    @{ Row = 17; D = '79.326.35';   E = '  +4.12%  ' }
    @{ Row = 18; D = '3.149.96';    E = '  +2.42%  ' }
    @{ Row = 19; D = '14.15';       E = '  +4.33%  ' }
    @{ Row = 20; D = '432.78';      E = '  +13.35%  ' }
    @{ Row = 21; D = '8.99';        E = '  +0.02%  ' }
    @{ Row = 22; D = '2.89';        E = '  +10.43%  ' }
    @{ Row = 23; D = '5.17';        E = '  +15.91%  ' }
    @{ Row = 24; D = '6.75';        E = '  +4.82%  ' }
    @{ Row = 25; D = '3.318.17';    E = '  +2.46%  ' }
    @{ Row = 26; D = '75.52';       E = '  +4.35%  ' }
    @{ Row = 27; D = '4.62';        E = '  +2.16%  ' }
    @{ Row = 28; D = '10.69';       E = '  +5.88%  ' }
    @{ Row = 29; D = '0.999';       E = '  -0.22%  ' }
    @{ Row = 30; D = '0.0000119';   E = '  +9.98%  ' }
    @{ Row = 31; D = '1.00';        E = '  +0.16%  ' }
    @{ Row = 32; D = '8.88';        E = '  +7.12%  ' }
    @{ Row = 33; D = '543.24';      E = '  +7.93%  ' }
    @{ Row = 34; D = '1.46';        E = '  +2.43%  ' }
    @{ Row = 35; D = '1.98';        E = '  +2.84%  ' }
    @{ Row = 36; D = '0.147';       E = '  +17.76%  ' }
    @{ Row = 37; D = '22.76';       E = '  +9.36%  ' }
    @{ Row = 38; D = '0.121';       E = '  +18.34%  ' }
    @{ Row = 39; D = '0.999';       E = '  -0.03%  ' }
    @{ Row = 40; D = '0.400';       E = '  +5.79%  ' }
    @{ Row = 41; D = '20.70';       E = '  +3.18%  ' }
    @{ Row = 42; D = '162.70';      E = '  +0.49%  ' }
    @{ Row = 43; D = '1.00';        E = '  -0.02%  ' }
    @{ Row = 44; D = '5.55';        E = '  +7.55%  ' }
    @{ Row = 45; D = '186.06';      E = '  -4.97%  ' }
    @{ Row = 46; D = '1.78';        E = '  +7.76%  ' }
    @{ Row = 47; D = '2.64';        E = '  +8.67%  ' }
    @{ Row = 48; D = '0.775';       E = '  -4.04%  ' }
    @{ Row = 49; D = '42.88';       E = '  +6.10%  ' }
    @{ Row = 50; D = '1.29';        E = '  +2.71%  ' }
    @{ Row = 51; D = '4.20';        E = '  +7.28%  ' }
)

foreach ($item in $data) {
    $dCell = $ws.Cells.Item($item.Row, 4)
    # Force text storage so numeric-looking prices (e.g. "1.00", "31.00")
    # keep their original formatting instead of being coerced to numbers.
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D

    $eCell = $ws.Cells.Item($item.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $item.E
}
